# edit.ps1 - applies the "updates files, adds privacy policy" commit:
#   1. Notes for slide 1: tweak a sentence ("might actually get used" -> "might get used")
#   2. Slide 6 "Title 1" placeholder: vertically center the text (anchor="ctr")
#   3. Slide 6 "Content Placeholder 2": vertically center the text (anchor="ctr")
#      and add a new bullet paragraph "Ask for feedback"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Notes slide 1 - shorten the sentence about the personal app download count.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$notes1 = $slide1.NotesPage
$notesShape = $notes1.Shapes.Item(2)
$notesRange = $notesShape.TextFrame.TextRange

# The notes placeholder has several paragraphs; PowerPoint reports them
# separated with CR (`r`), but only LF (`n`) round-trips correctly as a
# paragraph break when rewriting this particular text range, so normalize
# to `n` before writing the corrected sentence back.
$originalNotesText = $notesRange.Text
$notesParts = $originalNotesText -split "`r"
for ($i = 0; $i -lt $notesParts.Count; $i++) {
    if ($notesParts[$i] -eq "When a personal app of mine reached 30,000 downloads a few years back, I was really worried because I had no structure of maintaining that user base. And I realized that what I put out there into the wild, might actually get used by a lot of people.") {
        $notesParts[$i] = "When a personal app of mine reached 30,000 downloads a few years back, I was really worried because I had no structure of maintaining that user base. And I realized that what I put out there into the wild, might get used by a lot of people."
    }
}
$notesRange.Text = ($notesParts -join "`n")

# ---------------------------------------------------------------------------
# 2) & 3) Slide 6 ("Trust the team")
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)

# Title 1 -> vertical-center the text box.
$titleShape = $slide6.Shapes.Item(2)
$titleShape.TextFrame.VerticalAnchor = 3  # ppAnchorMiddle

# Content Placeholder 2 -> vertical-center + append a new bullet line.
$bodyShape = $slide6.Shapes.Item(3)
$bodyShape.TextFrame.VerticalAnchor = 3  # ppAnchorMiddle

$bodyRange = $bodyShape.TextFrame.TextRange
$lastParagraph = $bodyRange.Paragraphs($bodyRange.Paragraphs().Count, 1)
[void]$lastParagraph.InsertAfter("`rAsk for feedback")
